$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.420.97"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "'1.563.91"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("D4").Value = "'0.9999"

$ws.Range("D5").Value = "'0.9997"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").Value = "'285.88"
$ws.Range("E6").Value = "  -1.93%  "

$ws.Range("D7").Value = "'0.3646"
$ws.Range("E7").Value = "  -2.51%  "

$ws.Range("D8").Value = "'48.48"
$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("D9").Value = "'0.3340"
$ws.Range("E9").Value = "  -1.85%  "

$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").Value = "'0.07424"
$ws.Range("E11").Value = "  -1.75%  "

$ws.Range("D12").Value = "'1.0000"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").Value = "'20.85"
$ws.Range("E13").Value = "  -2.53%  "

$ws.Range("D14").Value = "'5.941"
$ws.Range("E14").Value = "  -0.90%  "

$ws.Range("D15").Value = "'6.899"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "'1.565.14"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "'0.00001107"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "'88.33"
$ws.Range("E18").Value = "  -3.03%  "

$ws.Range("D19").Value = "'0.06685"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("D21").Value = "'6.363"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("D22").Value = "'16.14"
$ws.Range("E22").Value = "  -1.66%  "

$ws.Range("D23").Value = "'12.00"
$ws.Range("E23").Value = "  -1.25%  "

$ws.Range("D24").Value = "'22.406.73"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").Value = "'2.414"
$ws.Range("E25").Value = "  +3.29%  "

$ws.Range("D26").Value = "'2.569"
$ws.Range("E26").Value = "  -1.06%  "

$ws.Range("D27").Value = "'149.98"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").Value = "'19.43"
$ws.Range("E28").Value = "  -3.48%  "

$ws.Range("D29").Value = "'5.006"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "'123.14"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").Value = "'1.738.93"
$ws.Range("E31").Value = "  -0.22%  "

$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("D33").Value = "'6.163"
$ws.Range("E33").Value = "  +0.66%  "

$ws.Range("D34").Value = "'1.991"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("D35").Value = "'9.793"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").Value = "'0.02403"
$ws.Range("E37").Value = "  -2.59%  "

$ws.Range("D38").Value = "'1.308"
$ws.Range("E38").Value = "  -5.64%  "

$ws.Range("D39").Value = "'0.06392"
$ws.Range("E39").Value = "  -2.29%  "

$ws.Range("D40").Value = "'0.2211"
$ws.Range("E40").Value = "  -3.48%  "

$ws.Range("D41").Value = "'5.344"
$ws.Range("E41").Value = "  -2.15%  "

$ws.Range("D42").Value = "'11.19"
$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").Value = "'13.80"
$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("D46").Value = "'3.761"
$ws.Range("E46").Value = "  -1.34%  "

$ws.Range("D47").Value = "'0.5771"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").Value = "'2.016"
$ws.Range("E48").Value = "  -3.45%  "

$ws.Range("D49").Value = "'125.02"
$ws.Range("E49").Value = "  -3.18%  "

$ws.Range("D50").Value = "'1.216"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("D51").Value = "'0.07215"
$ws.Range("E51").Value = "  -1.49%  "
